$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.733.21"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.528.64"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").Value = "2.920.08"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "2.554.34"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "42.810.83"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("D46").Value = "2.024.90"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "2.773.67"
$ws.Range("E51").Value = "  +0.46%  "
